# Auto-generated Excel COM-interop edit script
# Applies per-cell value updates to the "Shinryu Profits" workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the canonical-OOXML diff for the scheduled price-refresh commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates_ALC = @{
    "H64" = 3891.261
    "I64" = 3783.3333
    "J64" = 4009
    "K64" = 3783.3333
    "L64" = 4009
    "M64" = -3535.3333
    "N64" = -4505
    "H67" = 3891.261
    "I67" = 3783.3333
    "J67" = 4009
    "K67" = 3783.3333
    "L67" = 4009
    "M67" = -2925.3333
    "N67" = -5725
    "H76" = 3624.625
    "I76" = 3500
    "J76" = 3699.4
    "K76" = 3500
    "L76" = 3699.4
    "M76" = -3185
    "N76" = -4329.4
    "H79" = 3624.625
    "I79" = 3500
    "J79" = 3699.4
    "K79" = 3500
    "L79" = 3699.4
    "M79" = -2408
    "N79" = -5883.4
    "H125" = 3157.7273
    "I125" = 1999
    "J125" = 3273.6
    "K125" = 17991
    "L125" = 29462.4
    "M125" = -15531
    "N125" = -34382.39999999999
    "H129" = 1093.5566
    "I129" = 548.5
    "K129" = 1645.5
    "M129" = 3354.5
    "H135" = 891.1579
    "I135" = 672
    "J135" = 1713
    "K135" = 6048
    "L135" = 15417
    "M135" = -3513
    "N135" = -20487
}
foreach ($addr in $updates_ALC.Keys) {
    $ws.Range($addr).Value = $updates_ALC[$addr]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates_ARM = @{
    "H5" = 106.5
    "I5" = 74
    "J5" = 160.66667
    "K5" = 74
    "L5" = 160.66667
    "M5" = 38
    "N5" = -384.66667
    "H33" = 11013
    "I33" = 2026
    "J33" = 20000
    "K33" = 2026
    "L33" = 20000
    "M33" = -1697
    "N33" = -20658
    "H61" = 1839.88
    "I61" = 1263.45
    "J61" = 4145.6
    "K61" = 1263.45
    "L61" = 4145.6
    "M61" = -1051.45
    "N61" = -4569.6
    "H74" = 1426.9445
    "I74" = 1346.9678
    "J74" = 1922.8
    "K74" = 1346.9678
    "L74" = 1922.8
    "M74" = -472.9677999999999
    "N74" = -3670.8
    "H77" = 1426.9445
    "I77" = 1346.9678
    "J77" = 1922.8
    "K77" = 6734.839
    "L77" = 9614
    "M77" = -2366.839
    "N77" = -18350
    "H122" = 1854.091
    "I122" = 1882.7693
    "J122" = 1812.6666
    "K122" = 5648.3079
    "L122" = 5437.9998
    "M122" = -3198.3079
    "N122" = -10337.9998
    "H132" = 1949.5641
    "I132" = 1450.8928
    "J132" = 3218.9092
    "K132" = 4352.678400000001
    "L132" = 9656.7276
    "M132" = -1822.678400000001
    "N132" = -14716.7276
    "H136" = 1839.88
    "I136" = 1263.45
    "J136" = 4145.6
    "K136" = 3790.35
    "L136" = 12436.8
    "M136" = -1240.35
    "N136" = -17536.8
    "H139" = 45000
    "J139" = 45000
    "L139" = 45000
    "N139" = -55280
}
foreach ($addr in $updates_ARM.Keys) {
    $ws.Range($addr).Value = $updates_ARM[$addr]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates_BSM = @{
    "H4" = 106.5
    "I4" = 74
    "J4" = 160.66667
    "K4" = 74
    "L4" = 160.66667
    "M4" = 41
    "N4" = -390.66667
    "H22" = 103
    "I22" = 101
    "J22" = 105
    "K22" = 101
    "L22" = 105
    "M22" = 72
    "N22" = -451
    "H38" = 23344
    "I38" = 5032
    "J38" = 32500
    "K38" = 5032
    "L38" = 32500
    "M38" = -4616
    "N38" = -33332
    "H86" = 1821.2
    "I86" = 1821.2
    "J86" = 0
    "K86" = 1821.2
    "L86" = 0
    "M86" = -698.2
    "H89" = 1821.2
    "I89" = 1821.2
    "J89" = 0
    "K89" = 9106
    "L89" = 0
    "M89" = -3490
    "H112" = 38980
    "J112" = 38980
    "L112" = 38980
    "N112" = -41934
    "H134" = 2538.1765
    "I134" = 2384.3125
    "J134" = 5000
    "K134" = 7152.9375
    "L134" = 15000
    "M134" = -4617.9375
    "N134" = -20070
}
foreach ($addr in $updates_BSM.Keys) {
    $ws.Range($addr).Value = $updates_BSM[$addr]
}

$clears_BSM = @("N86", "N89")
foreach ($addr in $clears_BSM) {
    $ws.Range($addr).ClearContents()
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates_CRP = @{
    "H33" = 13985.134
    "I33" = 1092.75
    "K33" = 1092.75
    "M33" = -713.75
    "H112" = 35000
    "J112" = 35000
    "L112" = 35000
    "N112" = -37954
    "H141" = 31108.666
    "J141" = 30163
    "L141" = 30163
    "N141" = -40523
}
foreach ($addr in $updates_CRP.Keys) {
    $ws.Range($addr).Value = $updates_CRP[$addr]
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates_CUL = @{
    "H2" = 17.5
    "I2" = 9.625
    "J2" = 28
    "K2" = 57.75
    "L2" = 168
    "M2" = 55.25
    "N2" = -394
    "H3" = 4770
    "I3" = 3514
    "J3" = 5816.6665
    "K3" = 10542
    "L3" = 17449.9995
    "M3" = -10430
    "N3" = -17673.9995
    "H69" = 2603.9333
    "J69" = 2620
    "L69" = 7860
    "N69" = -9482
    "H72" = 2603.9333
    "J72" = 2620
    "L72" = 23580
    "N72" = -31692
    "H74" = 3675.3333
    "H77" = 3675.3333
    "H131" = 678.6429000000001
    "I131" = 413.5263
    "J131" = 897.65216
    "K131" = 1240.5789
    "L131" = 2692.95648
    "M131" = 3799.4211
    "N131" = -12772.95648
}
foreach ($addr in $updates_CUL.Keys) {
    $ws.Range($addr).Value = $updates_CUL[$addr]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates_GSM = @{
    "H102" = 1072.2273
    "I102" = 1066.0555
    "J102" = 1100
    "K102" = 1066.0555
    "L102" = 1100
    "M102" = 555.9445000000001
    "N102" = -4344
    "H113" = 7098
    "I113" = 1374.9
    "K113" = 1374.9
    "M113" = 795.0999999999999
    "H122" = 8334867
    "I122" = 12501176
    "J122" = 2250
    "K122" = 37503528
    "L122" = 6750
    "M122" = -37501078
    "N122" = -11650
    "H126" = 4007.5557
    "I126" = 2753
    "K126" = 8259
    "M126" = -5789
    "H38" = 5981.5
    "I38" = 4030
    "J38" = 7933
    "K38" = 4030
    "L38" = 7933
    "M38" = -3620
    "N38" = -8753
}
foreach ($addr in $updates_GSM.Keys) {
    $ws.Range($addr).Value = $updates_GSM[$addr]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates_LTW = @{
    "H40" = 6277.5386
    "I40" = 5856.4443
    "J40" = 7225
    "K40" = 5856.4443
    "L40" = 7225
    "M40" = -5720.4443
    "N40" = -7497
    "H132" = 2592.1738
    "I132" = 1697.8572
    "J132" = 5437.727
    "K132" = 5093.571599999999
    "L132" = 16313.181
    "M132" = -2563.571599999999
    "N132" = -21373.181
    "H136" = 2032.8966
    "I136" = 1663.762
    "J136" = 3001.875
    "K136" = 4991.286
    "L136" = 9005.625
    "M136" = -2441.286
    "N136" = -14105.625
    "H140" = 0
    "J140" = 0
    "L140" = 0
    "H141" = 0
    "J141" = 0
    "L141" = 0
}
foreach ($addr in $updates_LTW.Keys) {
    $ws.Range($addr).Value = $updates_LTW[$addr]
}

$clears_LTW = @("N140", "N141")
foreach ($addr in $clears_LTW) {
    $ws.Range($addr).ClearContents()
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates_WVR = @{
    "H28" = 70019
    "I28" = 0
    "J28" = 70019
    "K28" = 0
    "L28" = 70019
    "N28" = -70715
    "H122" = 2925.6875
    "I122" = 2965.6667
    "J122" = 2874.2856
    "K122" = 8897.000100000001
    "L122" = 8622.856800000001
    "M122" = -6447.000100000001
    "N122" = -13522.8568
    "H133" = 0
    "J133" = 0
    "L133" = 0
}
foreach ($addr in $updates_WVR.Keys) {
    $ws.Range($addr).Value = $updates_WVR[$addr]
}

$clears_WVR = @("M28", "N133")
foreach ($addr in $clears_WVR) {
    $ws.Range($addr).ClearContents()
}
